# "hierarchisation des balises de titres"
#
# 1) Remove the bullet "Pas de h2 directement h3" entirely.
# 2) Swap the order of the two bullets "Rajouter ALT(12/16) et Title
#    (16/16)dans image" and "Ajouter title dans les liens" so that the
#    "liens" bullet now comes first.

$d = $word.ActiveDocument

# --- Step 1: delete the "Pas de h2 directement h3" paragraph -------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "Pas de h2 directement h3") {
        $p.Range.Delete()
        break
    }
}

# --- Step 2: swap "Rajouter ALT(...)...image" and "Ajouter title dans
#             les liens" bullets ------------------------------------------
$imageParagraph = $null
$linksParagraph = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "Rajouter ALT(12/16) et Title (16/16)dans image") {
        $imageParagraph = $p
    }
    elseif ($t -eq "Ajouter title dans les liens") {
        $linksParagraph = $p
    }
}

if ($imageParagraph -ne $null -and $linksParagraph -ne $null) {
    $imageRange = $d.Range($imageParagraph.Range.Start, $imageParagraph.Range.End - 1)
    $linksRange = $d.Range($linksParagraph.Range.Start, $linksParagraph.Range.End - 1)

    $imageText = $imageRange.Text
    $linksText = $linksRange.Text

    $linksRange.Text = $imageText
    $imageRange.Text = $linksText
}
